# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status for the c143638d row flips from "Ready for handoff" to
# "Handback transform failed" everywhere it is reported: the per-locale
# status columns on the Overview sheet and the Status column on each
# locale's own sheet.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# zh-cn: Error Detail for the c143638d row (row 3) gets a transform-mismatch message
$zhcn.Range("P3").Value = "Handback file name: uvcwu0j2.xkc is different with handoff file name: c143638d-0e11-4537-9322-4e6a76dc6aa4.2f47a77239c52f17f661fdbed264660ff4c7cf24.zh-cn."
# Stored column width of 40 chars (COM ColumnWidth snaps to a pixel grid, so
# 235/6 is the value that round-trips to a stored width of exactly 40)
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# de-de: Error Detail for the c143638d row (row 3) gets a transform-mismatch message
$dede.Range("P3").Value = "Handback file name: uvcwu0j2.xkc is different with handoff file name: c143638d-0e11-4537-9322-4e6a76dc6aa4.2f47a77239c52f17f661fdbed264660ff4c7cf24.de-de."
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
